$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph near the top of the document
#    <w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#         <w:r><w:t>: Read our review ...</w:t></w:r></w:p>
# Locate the paragraph containing "Meta description" and delete the whole paragraph (incl. mark)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description:*") {
        $p.Range.Delete()
        break
    }
}

# 2. At the end of the document, insert a new bold paragraph "Play Disco Diamonds Free
#    Slot: Review & Similar Games" right before the final "Prompt: ..." paragraph, and
#    rewrite that final paragraph's text while keeping its italic run formatting.
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)

$last.Range.InsertParagraphBefore() | Out-Null

$count2 = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count2 - 1)

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Disco Diamonds Free Slot: Review &amp; Similar Games</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($xmlFrag) | Out-Null

$oldText = 'Prompt: Create a cartoon-style feature image that showcases the fun and unique identity of the "Disco Diamonds" game. The image should focus on a happy Maya warrior character wearing glasses, as a nod to the game''s disco party theme. The image should be bright and colorful with funky disco elements, such as a disco ball and neon lights. The Maya warrior character should be drawn with a big smile and wearing stylish glasses, emphasizing the fun, upbeat feel of the game. The image should also include text that reads "Disco Diamonds", using a bold and playful font. The text should stand out and capture the attention of potential players. Overall, the feature image should capture the essence of the game''s party theme while also showcasing its simple and engaging gameplay. The colorful and eye-catching design should entice players to give the game a try and experience the excitement of Disco Diamonds.'
$newText = 'Read our review of Disco Diamonds and play it for free. Find similar party slots like Disco Danny and Dance Party.'

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
